# "Drafting several 'red flag' measures"
# Mark a handful of measures as "Draft" in the Done column (K):
#   - Contracts not notified      (row 2)
#   - Median price change         (row 3)
#   - Mean price change           (row 4)
#   - Purchase specificity (count) (row 19)
#   - Purchase specificity (value) (row 20)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measures")

$ws.Range("K2").Value = "Draft"
$ws.Range("K3").Value = "Draft"
$ws.Range("K4").Value = "Draft"
$ws.Range("K19").Value = "Draft"
$ws.Range("K20").Value = "Draft"

# Restore the frozen-pane view (row 1 / column A) and update the
# scroll/selection state to match where the author left the cursor.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$ws.Range("B16").Select()
